$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix row 9 (S.No. 6) values and add missing Boundary / Expected Output labels ---
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "Min"
$ws.Range("F9").Value = "equilateral"

# --- Row 10 (S.No. 7) - new boundary test data ---
$ws.Range("B10").Value = -50
$ws.Range("C10").Value = -50
$ws.Range("D10").Value = -50
$ws.Range("E10").Value = "Min -"
$ws.Range("F10").Value = "invalid"

# --- Row 11 (S.No. 8) - new boundary test data ---
$ws.Range("E11").Value = "Max"
$ws.Range("B11").Value = "Int32.MaxValue"
$ws.Range("C11").Value = "Int32.MaxValue"
$ws.Range("D11").Value = "Int32.MaxValue"
$ws.Range("F11").Value = "invalid"

# --- Column width adjustments: col A stays 11, cols B:D widen to fit new (longer) content ---
$ws.Range("B:D").ColumnWidth = 14.14

# --- Move the active selection down one row, to reflect the newly filled row ---
$ws.Range("F14").Select() | Out-Null
